# Median edu by county type
# Adds a "Median pct BA" column (K) to the Summary Sheet, with header,
# per-row data values, a light-purple color-scale conditional format,
# and adjusts dimension / selection / column widths / merged header cell
# to extend from column J to column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary Sheet")

# --- Header cell K2: copy J2's header format, then add wrap text (matches
# the "center + wrapText" look used for the new header) and set the text.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K2").WrapText = $true
$ws.Range("K2").Value = "Median pct BA"

# --- K1 is part of the merged title row, same style as the rest of row 1.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data cells K3:K16 / K17: plain bordered cells (same border family as
# column J), holding the median pct-BA value for each county type.
$ws.Range("K3").Style = "Normal"
$ws.Range("K3").Borders.Item(7).LineStyle = 1
$ws.Range("K3").Borders.Item(7).Weight = 2
$ws.Range("K3").Borders.Item(10).LineStyle = 1
$ws.Range("K3").Borders.Item(10).Weight = 2
$ws.Range("K4:K16").Style = "Normal"
$ws.Range("K4:K16").Borders.Item(7).LineStyle = 1
$ws.Range("K4:K16").Borders.Item(7).Weight = 2
$ws.Range("K4:K16").Borders.Item(10).LineStyle = 1
$ws.Range("K4:K16").Borders.Item(10).Weight = 2
$ws.Range("K17").Style = "Normal"
$ws.Range("K17").Borders.Item(7).LineStyle = 1
$ws.Range("K17").Borders.Item(7).Weight = 2
$ws.Range("K17").Borders.Item(10).LineStyle = -4119
$ws.Range("K17").Borders.Item(10).Weight = -4138

$kValues = @(0.28, 0.201, 0.144, 0.151, 0.145, 0.2415, 0.357, 0.15, 0.1595, 0.198, 0.334, 0.212, 0.201, 0.323, 0.251)
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
}

# --- Column widths: J and K both become 12.140625 wide (was 10 for J only).
$ws.Range("J1:K1").ColumnWidth = 12.140625

# --- Extend the merged title cell and the print selection from J to K.
$ws.Range("B1:J1").UnMerge()
$ws.Range("B1:K1").Merge()
$ws.Range("B1:K1").Select()

# --- Color-scale conditional formatting on the new column (light purple).
$kRange = $ws.Range("K3:K17")
$cf = $kRange.FormatConditions.AddColorScale(2)
$cf.ColorScaleCriteria.Item(1).Type = 1  # xlConditionValueLowestValue
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 16777212   # FFFCFCFF (BGR)
$cf.ColorScaleCriteria.Item(2).Type = 2  # xlConditionValueHighestValue
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 16750233   # FF9966FF (BGR)

Write-Host "Median pct BA column added"
